$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($addr, $val) {
    $c = $ws.Range($addr)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.Style = "Normal"
}

Set-TextCell "D2" "62.347.52"
$ws.Range("E2").Value = "  -1.66%  "
Set-TextCell "D3" "2.454.19"
$ws.Range("E3").Value = "  -0.20%  "
$ws.Range("E4").Value = "  -0.08%  "
Set-TextCell "D5" "581.09"
$ws.Range("E5").Value = "  +1.33%  "
Set-TextCell "D6" "143.64"
$ws.Range("E6").Value = "  -2.77%  "
$ws.Range("E7").Value = "  -0.04%  "
Set-TextCell "D8" "0.530"
$ws.Range("E8").Value = "  -1.66%  "
Set-TextCell "D9" "2.453.75"
$ws.Range("E9").Value = "  -0.29%  "
$ws.Range("E10").Value = "  -3.47%  "
Set-TextCell "D11" "0.159"
$ws.Range("E11").Value = "  +1.55%  "
Set-TextCell "D12" "5.21"
$ws.Range("E12").Value = "  -1.07%  "
Set-TextCell "D13" "0.346"
$ws.Range("E13").Value = "  -2.76%  "
Set-TextCell "D14" "26.48"
$ws.Range("E14").Value = "  -2.51%  "
Set-TextCell "D15" "0.0000173"
$ws.Range("E15").Value = "  -4.18%  "
Set-TextCell "D16" "2.865.00"
$ws.Range("E16").Value = "  -1.50%  "
Set-TextCell "D17" "62.309.80"
$ws.Range("E17").Value = "  -1.72%  "
Set-TextCell "D18" "2.448.53"
$ws.Range("E18").Value = "  -0.56%  "
Set-TextCell "D19" "10.98"
$ws.Range("E19").Value = "  -3.52%  "
Set-TextCell "D20" "7.14"
$ws.Range("E20").Value = "  -3.28%  "
Set-TextCell "D21" "330.87"
$ws.Range("E21").Value = "  +0.35%  "
Set-TextCell "D22" "4.13"
$ws.Range("E22").Value = "  -1.81%  "
Set-TextCell "D23" "1.96"
$ws.Range("E23").Value = "  -5.55%  "
$ws.Range("E24").Value = "  +0.11%  "
Set-TextCell "D25" "65.71"
$ws.Range("E25").Value = "  -0.03%  "
Set-TextCell "D26" "9.44"
$ws.Range("E26").Value = "  +4.17%  "
Set-TextCell "D27" "628.14"
$ws.Range("E27").Value = "  +0.30%  "
Set-TextCell "D29" "0.0₃0953"
$ws.Range("E29").Value = "  -8.31%  "
$ws.Range("E30").Value = "  -0.29%  "
$ws.Range("E31").Value = "  -6.12%  "
Set-TextCell "D32" "8.03"
$ws.Range("E32").Value = "  -3.74%  "
$ws.Range("E33").Value = "  +0.14%  "
$ws.Range("E34").Value = "  -1.16%  "
Set-TextCell "D35" "4.95"
$ws.Range("E35").Value = "  -5.26%  "
$ws.Range("E36").Value = "  +0.22%  "
Set-TextCell "D37" "1.44"
$ws.Range("E37").Value = "  -6.99%  "
Set-TextCell "D38" "0.376"
$ws.Range("E38").Value = "  -1.35%  "
Set-TextCell "D39" "150.34"
$ws.Range("E39").Value = "  +3.44%  "
Set-TextCell "D40" "18.37"
$ws.Range("E40").Value = "  -2.40%  "
Set-TextCell "D41" "5.27"
$ws.Range("E41").Value = "  -4.00%  "
Set-TextCell "D42" "1.77"
$ws.Range("E42").Value = "  -1.91%  "
Set-TextCell "D43" "42.78"
$ws.Range("E43").Value = "  +2.18%  "
$ws.Range("E44").Value = "  +0.00%  "
Set-TextCell "D45" "2.49"
$ws.Range("E45").Value = "  -8.07%  "
Set-TextCell "D46" "143.45"
$ws.Range("E46").Value = "  -4.32%  "
Set-TextCell "D47" "3.65"
$ws.Range("E47").Value = "  -3.27%  "
Set-TextCell "D48" "0.0525"
$ws.Range("E48").Value = "  -2.41%  "
Set-TextCell "D49" "0.602"
$ws.Range("E49").Value = "  -0.10%  "
Set-TextCell "D50" "19.64"
$ws.Range("E50").Value = "  -8.29%  "
Set-TextCell "D51" "0.0₆0235"
$ws.Range("E51").Value = "  +4.44%  "
